$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 437.33334
$ws.Range("I12").Value = 108
$ws.Range("K12").Value = 108
$ws.Range("M12").Value = 62
$ws.Range("H100").Value = 1893.1111
$ws.Range("I100").Value = 1552.1538
$ws.Range("K100").Value = 1552.1538
$ws.Range("M100").Value = -1011.1538
$ws.Range("H129").Value = 1737.909
$ws.Range("J129").Value = 2553.6
$ws.Range("L129").Value = 7660.799999999999
$ws.Range("N129").Value = -17660.8
$ws.Range("H137").Value = 1789.909
$ws.Range("I137").Value = 959.46155
$ws.Range("K137").Value = 2878.38465
$ws.Range("M137").Value = -328.38465

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2299.6667
$ws.Range("I2").Value = 2299.6667
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 2299.6667
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -2186.6667
$ws.Range("N2").ClearContents()
$ws.Range("H98").Value = 6188.25
$ws.Range("I98").Value = 1000
$ws.Range("J98").Value = 7917.6665
$ws.Range("K98").Value = 1000
$ws.Range("L98").Value = 7917.6665
$ws.Range("M98").Value = 1995
$ws.Range("N98").Value = -13907.6665
$ws.Range("H116").Value = 2299.6667
$ws.Range("I116").Value = 2299.6667
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 2299.6667
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = -5.666700000000219
$ws.Range("N116").ClearContents()

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2299.6667
$ws.Range("I3").Value = 2299.6667
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 2299.6667
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -2185.6667
$ws.Range("N3").ClearContents()
$ws.Range("H20").Value = 3600.7144
$ws.Range("I20").Value = 4083
$ws.Range("J20").Value = 707
$ws.Range("K20").Value = 4083
$ws.Range("L20").Value = 707
$ws.Range("M20").Value = -3836
$ws.Range("N20").Value = -1201
$ws.Range("H22").Value = 437.6
$ws.Range("I22").Value = 429.66666
$ws.Range("J22").Value = 449.5
$ws.Range("K22").Value = 429.66666
$ws.Range("L22").Value = 449.5
$ws.Range("M22").Value = -256.66666
$ws.Range("N22").Value = -795.5
$ws.Range("H86").Value = 4007
$ws.Range("I86").Value = 1523.625
$ws.Range("K86").Value = 1523.625
$ws.Range("M86").Value = -400.625
$ws.Range("H89").Value = 4007
$ws.Range("I89").Value = 1523.625
$ws.Range("K89").Value = 7618.125
$ws.Range("M89").Value = -2002.125
$ws.Range("H106").Value = 17633.334
$ws.Range("J106").Value = 17633.334
$ws.Range("L106").Value = 17633.334
$ws.Range("N106").Value = -20157.334
$ws.Range("H137").Value = 124500
$ws.Range("J137").Value = 124500
$ws.Range("L137").Value = 124500
$ws.Range("N137").Value = -134700

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 73.46666999999999
$ws.Range("I7").Value = 50.666668
$ws.Range("K7").Value = 50.666668
$ws.Range("M7").Value = 62.333332
$ws.Range("H16").Value = 1183.2858
$ws.Range("I16").Value = 1183.2858
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 1183.2858
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -896.2858000000001
$ws.Range("N16").ClearContents()
$ws.Range("H22").Value = 2962.7144
$ws.Range("I22").Value = 1622.25
$ws.Range("K22").Value = 1622.25
$ws.Range("M22").Value = -1272.25
$ws.Range("H58").Value = 1274.6666
$ws.Range("I58").Value = 1347.4
$ws.Range("K58").Value = 1347.4
$ws.Range("M58").Value = -1144.4
$ws.Range("H62").Value = 2737
$ws.Range("I62").Value = 2737
$ws.Range("K62").Value = 2737
$ws.Range("M62").Value = -2113
$ws.Range("H65").Value = 2737
$ws.Range("I65").Value = 2737
$ws.Range("K65").Value = 13685
$ws.Range("M65").Value = -10565
$ws.Range("H86").Value = 5500
$ws.Range("I86").Value = 5000
$ws.Range("J86").Value = 6000
$ws.Range("K86").Value = 5000
$ws.Range("L86").Value = 6000
$ws.Range("M86").Value = -3877
$ws.Range("N86").Value = -8246
$ws.Range("H89").Value = 5500
$ws.Range("I89").Value = 5000
$ws.Range("J89").Value = 6000
$ws.Range("K89").Value = 25000
$ws.Range("L89").Value = 30000
$ws.Range("M89").Value = -19384
$ws.Range("N89").Value = -41232
$ws.Range("H94").Value = 3648.3635
$ws.Range("I94").Value = 1019.8
$ws.Range("K94").Value = 1019.8
$ws.Range("M94").Value = -568.8
$ws.Range("H113").Value = 1183.2858
$ws.Range("I113").Value = 1183.2858
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1183.2858
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 986.7141999999999
$ws.Range("N113").ClearContents()
$ws.Range("H136").Value = 1274.6666
$ws.Range("I136").Value = 1347.4
$ws.Range("K136").Value = 4042.2
$ws.Range("M136").Value = -1492.2

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H57").Value = 1999.5
$ws.Range("I57").Value = 1999.5
$ws.Range("J57").Value = 0
$ws.Range("K57").Value = 5998.5
$ws.Range("L57").Value = 0
$ws.Range("M57").Value = -5439.5
$ws.Range("N57").ClearContents()
$ws.Range("H113").Value = 1482
$ws.Range("J113").Value = 1719.4445
$ws.Range("L113").Value = 5158.333500000001
$ws.Range("N113").Value = -9498.333500000001

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("M70").ClearContents()
$ws.Range("N70").ClearContents()
$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("M73").ClearContents()
$ws.Range("N73").ClearContents()
$ws.Range("H98").Value = 7797.8184
$ws.Range("I98").Value = 9000
$ws.Range("J98").Value = 7677.6
$ws.Range("K98").Value = 9000
$ws.Range("L98").Value = 7677.6
$ws.Range("M98").Value = -6005
$ws.Range("N98").Value = -13667.6
$ws.Range("H99").Value = 2050.5715
$ws.Range("I99").Value = 2050.5715
$ws.Range("K99").Value = 2050.5715
$ws.Range("M99").Value = 195.4285
$ws.Range("H102").Value = 2209.5715
$ws.Range("I102").Value = 2209.5715
$ws.Range("K102").Value = 2209.5715
$ws.Range("M102").Value = -587.5715
$ws.Range("H122").Value = 1900.5625
$ws.Range("I122").Value = 1877.6923
$ws.Range("K122").Value = 5633.0769
$ws.Range("M122").Value = -3183.0769
$ws.Range("H126").Value = 2921.2144
$ws.Range("I126").Value = 2921.2144
$ws.Range("K126").Value = 8763.643199999999
$ws.Range("M126").Value = -6293.643199999999

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3468.625
$ws.Range("I7").Value = 3233.1667
$ws.Range("J7").Value = 4175
$ws.Range("K7").Value = 3233.1667
$ws.Range("L7").Value = 4175
$ws.Range("M7").Value = -3121.1667
$ws.Range("N7").Value = -4399
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("M27").ClearContents()
$ws.Range("H40").Value = 5300.4
$ws.Range("I40").Value = 4626.75
$ws.Range("K40").Value = 4626.75
$ws.Range("M40").Value = -4490.75
$ws.Range("H46").Value = 2805.524
$ws.Range("J46").Value = 3595.077
$ws.Range("L46").Value = 3595.077
$ws.Range("N46").Value = -3971.077
$ws.Range("H55").Value = 1180.6428
$ws.Range("I55").Value = 1426.5555
$ws.Range("J55").Value = 738
$ws.Range("K55").Value = 1426.5555
$ws.Range("L55").Value = 738
$ws.Range("M55").Value = -1253.5555
$ws.Range("N55").Value = -1084
$ws.Range("H68").Value = 7061.5
$ws.Range("I68").Value = 4123
$ws.Range("J68").Value = 10000
$ws.Range("K68").Value = 4123
$ws.Range("L68").Value = 10000
$ws.Range("M68").Value = -3374
$ws.Range("N68").Value = -11498
$ws.Range("H71").Value = 7061.5
$ws.Range("I71").Value = 4123
$ws.Range("J71").Value = 10000
$ws.Range("K71").Value = 20615
$ws.Range("L71").Value = 50000
$ws.Range("M71").Value = -16871
$ws.Range("N71").Value = -57488
$ws.Range("H93").Value = 1019.44
$ws.Range("I93").Value = 956.2381
$ws.Range("K93").Value = 956.2381
$ws.Range("M93").Value = 291.7619
$ws.Range("H126").Value = 3468.625
$ws.Range("I126").Value = 3233.1667
$ws.Range("J126").Value = 4175
$ws.Range("K126").Value = 9699.500100000001
$ws.Range("L126").Value = 12525
$ws.Range("M126").Value = -7229.500100000001
$ws.Range("N126").Value = -17465
$ws.Range("H132").Value = 9842.857
$ws.Range("J132").Value = 17502.5
$ws.Range("L132").Value = 52507.5
$ws.Range("N132").Value = -57567.5

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 999.6667
$ws.Range("I81").Value = 999.6667
$ws.Range("K81").Value = 1999.3334
$ws.Range("M81").Value = -938.3334
$ws.Range("H84").Value = 999.6667
$ws.Range("I84").Value = 999.6667
$ws.Range("K84").Value = 9996.666999999999
$ws.Range("M84").Value = -4692.666999999999
$ws.Range("H126").Value = 4581.2
$ws.Range("I126").Value = 2369.6667
$ws.Range("K126").Value = 7109.000100000001
$ws.Range("M126").Value = -4639.000100000001
